$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
    $ws.Range($addr).Style = "Normal"
}

function Set-Cell($addr, $value) {
    $ws.Range($addr).Value = $value
}

# Row 2
Set-Cell "D2" "45.550.83"
Set-Cell "E2" "  +6.71%  "

# Row 3
Set-Cell "D3" "2.400.55"
Set-Cell "E3" "  +4.25%  "

# Row 4
Set-Cell "E4" "  +0.07%  "

# Row 5
Set-TextCell "D5" "114.20"
Set-Cell "E5" "  +8.67%  "

# Row 6
Set-TextCell "D6" "320.96"
Set-Cell "E6" "  +3.61%  "

# Row 7
Set-TextCell "D7" "0.635"
Set-Cell "E7" "  +2.53%  "

# Row 8
Set-Cell "E8" "  -0.07%  "

# Row 9
Set-Cell "E9" "  +3.44%  "

# Row 10
Set-TextCell "D10" "42.43"
Set-Cell "E10" "  +6.90%  "

# Row 11
Set-Cell "E11" "  +3.17%  "

# Row 12
Set-TextCell "D12" "8.68"
Set-Cell "E12" "  +4.99%  "

# Row 13
Set-Cell "E13" "  +2.96%  "

# Row 14
Set-Cell "E14" "  +1.72%  "

# Row 15
Set-TextCell "D15" "15.95"
Set-Cell "E15" "  +4.36%  "

# Row 16
Set-Cell "D16" "2.762.40"
Set-Cell "E16" "  -0.74%  "

# Row 17
Set-Cell "D17" "2.400.68"
Set-Cell "E17" "  +4.61%  "

# Row 18
Set-Cell "D18" "45.495.13"
Set-Cell "E18" "  +6.66%  "

# Row 19
Set-TextCell "D19" "7.53"
Set-Cell "E19" "  +2.91%  "

# Row 20
Set-Cell "E20" "  +3.25%  "

# Row 21
Set-TextCell "D21" "13.47"
Set-Cell "E21" "  -1.45%  "

# Row 22
Set-TextCell "D22" "74.88"
Set-Cell "E22" "  +2.07%  "

# Row 23
Set-TextCell "D23" "3.59"
Set-Cell "E23" "  +4.62%  "

# Row 24
Set-TextCell "D24" "264.66"
Set-Cell "E24" "  -0.69%  "

# Row 25
Set-TextCell "D25" "2.36"
Set-Cell "E25" "  +5.34%  "

# Row 26
Set-Cell "E26" "  -0.53%  "

# Row 27
Set-TextCell "D27" "7.69"
Set-Cell "E27" "  +5.29%  "

# Row 28
Set-TextCell "D28" "11.38"
Set-Cell "E28" "  +3.49%  "

# Row 29
Set-Cell "E29" "  +3.07%  "

# Row 30
Set-TextCell "D30" "39.97"
Set-Cell "E30" "  +8.40%  "

# Row 31
Set-TextCell "D31" "22.79"
Set-Cell "E31" "  +2.45%  "

# Row 32
Set-TextCell "D32" "0.0976"
Set-Cell "E32" "  +13.63%  "

# Row 33
Set-TextCell "D33" "173.13"
Set-Cell "E33" "  +5.22%  "

# Row 34
Set-TextCell "D34" "2.94"
Set-Cell "E34" "  +3.08%  "

# Row 35
Set-TextCell "D35" "0.133"
Set-Cell "E35" "  +1.61%  "

# Row 36
Set-TextCell "D36" "4.93"
Set-Cell "E36" "  +7.96%  "

# Row 37
Set-Cell "E37" "  +6.45%  "

# Row 38
Set-TextCell "D38" "4.12"
Set-Cell "E38" "  +13.66%  "

# Row 39
Set-TextCell "D39" "3.08"
Set-Cell "E39" "  +8.32%  "

# Row 40
Set-Cell "E40" "  +4.51%  "

# Row 41
Set-TextCell "D41" "1.79"
Set-Cell "E41" "  +13.76%  "

# Row 42
Set-Cell "B42" "Celestia"
Set-Cell "C42" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextCell "D42" "13.83"
Set-Cell "E42" "  +12.89%  "

# Row 43
Set-Cell "B43" "Algorand"
Set-Cell "C43" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell "D43" "0.241"
Set-Cell "E43" "  +5.53%  "

# Row 44
Set-TextCell "D44" "100.27"
Set-Cell "E44" "  -7.75%  "

# Row 45
Set-TextCell "D45" "71.58"
Set-Cell "E45" "  -1.01%  "

# Row 46
Set-TextCell "D46" "88.05"
Set-Cell "E46" "  +15.04%  "

# Row 47
Set-Cell "E47" "  +0.02%  "

# Row 48
Set-TextCell "D48" "5.82"
Set-Cell "E48" "  +13.17%  "

# Row 49
Set-TextCell "D49" "116.34"
Set-Cell "E49" "  +4.61%  "

# Row 50
Set-TextCell "D50" "9.51"
Set-Cell "E50" "  +9.34%  "

# Row 51
Set-Cell "B51" "Maker"
Set-Cell "C51" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-Cell "D51" "1.662.48"
Set-Cell "E51" "  -3.53%  "
